$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2800.1428
$ws.Range("I18").Value = 2850.1667
$ws.Range("J18").Value = 2500
$ws.Range("K18").Value = 2850.1667
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = -2566.1667
$ws.Range("N18").Value = -3068

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1862.5
$ws.Range("I32").Value = 1200.4
$ws.Range("J32").Value = 2335.4285
$ws.Range("K32").Value = 1200.4
$ws.Range("L32").Value = 2335.4285
$ws.Range("M32").Value = -874.4000000000001
$ws.Range("N32").Value = -2987.4285

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 129.17392
$ws.Range("I33").Value = 86.42104999999999
$ws.Range("J33").Value = 332.25
$ws.Range("K33").Value = 86.42104999999999
$ws.Range("L33").Value = 332.25
$ws.Range("M33").Value = 142.57895
$ws.Range("N33").Value = -790.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 192.6
$ws.Range("I55").Value = 165.14285
$ws.Range("J55").Value = 256.66666
$ws.Range("K55").Value = 165.14285
$ws.Range("L55").Value = 256.66666
$ws.Range("M55").Value = 48.85714999999999
$ws.Range("N55").Value = -684.66666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1083.3334
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 1125
$ws.Range("K125").Value = 9000
$ws.Range("L125").Value = 10125
$ws.Range("M125").Value = -6540
$ws.Range("N125").Value = -15045

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2021.2646
$ws.Range("I137").Value = 1125.5625
$ws.Range("J137").Value = 2817.4443
$ws.Range("K137").Value = 3376.6875
$ws.Range("L137").Value = 8452.332900000001
$ws.Range("M137").Value = -826.6875
$ws.Range("N137").Value = -13552.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1836.72
$ws.Range("I45").Value = 1658
$ws.Range("J45").Value = 2551.6
$ws.Range("K45").Value = 1658
$ws.Range("L45").Value = 2551.6
$ws.Range("M45").Value = -1281

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 993.44446
$ws.Range("I74").Value = 849.7727
$ws.Range("J74").Value = 1625.6
$ws.Range("K74").Value = 849.7727
$ws.Range("L74").Value = 1625.6
$ws.Range("M74").Value = 24.22730000000001
$ws.Range("N74").Value = -3373.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 993.44446
$ws.Range("I77").Value = 849.7727
$ws.Range("J77").Value = 1625.6
$ws.Range("K77").Value = 4248.863499999999
$ws.Range("L77").Value = 8128
$ws.Range("M77").Value = 119.1365000000005
$ws.Range("N77").Value = -16864

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2006
$ws.Range("I122").Value = 2006
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6018
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3568
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2355.7778
$ws.Range("I132").Value = 2130.5278
$ws.Range("J132").Value = 3256.7778
$ws.Range("K132").Value = 6391.5834
$ws.Range("L132").Value = 9770.3334
$ws.Range("M132").Value = -3861.5834
$ws.Range("N132").Value = -14830.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 146.625
$ws.Range("I22").Value = 139
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 139
$ws.Range("L22").Value = 200
$ws.Range("M22").Value = 34

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2506.3696
$ws.Range("I86").Value = 2961.3704
$ws.Range("J86").Value = 1859.7894
$ws.Range("K86").Value = 2961.3704
$ws.Range("L86").Value = 1859.7894
$ws.Range("M86").Value = -1838.3704
$ws.Range("N86").Value = -4105.7894

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2506.3696
$ws.Range("I89").Value = 2961.3704
$ws.Range("J89").Value = 1859.7894
$ws.Range("K89").Value = 14806.852
$ws.Range("L89").Value = 9298.947
$ws.Range("M89").Value = -9190.851999999999
$ws.Range("N89").Value = -20530.947

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 865.06665
$ws.Range("I31").Value = 641.6977000000001
$ws.Range("J31").Value = 1430.0588
$ws.Range("K31").Value = 641.6977000000001
$ws.Range("L31").Value = 1430.0588
$ws.Range("M31").Value = -346.6977000000001
$ws.Range("N31").Value = -2020.0588

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 865.06665
$ws.Range("I34").Value = 641.6977000000001
$ws.Range("J34").Value = 1430.0588
$ws.Range("K34").Value = 641.6977000000001
$ws.Range("L34").Value = 1430.0588
$ws.Range("M34").Value = -439.6977000000001
$ws.Range("N34").Value = -1834.0588

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 7826.3335
$ws.Range("I41").Value = 3991.6
$ws.Range("J41").Value = 27000
$ws.Range("K41").Value = 3991.6
$ws.Range("L41").Value = 27000
$ws.Range("M41").Value = -3563.6
$ws.Range("N41").Value = -27856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 904.5
$ws.Range("I58").Value = 698.13336
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 698.13336
$ws.Range("L58").Value = 4000
$ws.Range("M58").Value = -495.13336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 22500
$ws.Range("I59").Value = 15000
$ws.Range("J59").Value = 25000
$ws.Range("K59").Value = 15000
$ws.Range("L59").Value = 25000
$ws.Range("M59").Value = -13855
$ws.Range("N59").Value = -27290

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 22714.285
$ws.Range("I60").Value = 10000
$ws.Range("J60").Value = 24833.334
$ws.Range("K60").Value = 10000
$ws.Range("L60").Value = 24833.334
$ws.Range("M60").Value = -9489
$ws.Range("N60").Value = -25855.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 500
$ws.Range("I122").Value = 500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 950

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 11112345
$ws.Range("I134").Value = 13334394
$ws.Range("J134").Value = 2100
$ws.Range("K134").Value = 40003182
$ws.Range("L134").Value = 6300
$ws.Range("M134").Value = -40000647

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 904.5
$ws.Range("I136").Value = 698.13336
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 2094.40008
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = 455.5999199999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 660.7692
$ws.Range("I18").Value = 644
$ws.Range("J18").Value = 716.6667
$ws.Range("K18").Value = 1932
$ws.Range("L18").Value = 2150.0001
$ws.Range("M18").Value = -1763
$ws.Range("N18").Value = -2488.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2298.5264
$ws.Range("I68").Value = 900
$ws.Range("J68").Value = 2463.0588
$ws.Range("K68").Value = 2700
$ws.Range("L68").Value = 7389.176399999999
$ws.Range("M68").Value = -1889
$ws.Range("N68").Value = -9011.1764

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2298.5264
$ws.Range("I71").Value = 900
$ws.Range("J71").Value = 2463.0588
$ws.Range("K71").Value = 8100
$ws.Range("L71").Value = 22167.5292
$ws.Range("M71").Value = -4044
$ws.Range("N71").Value = -30279.5292

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1000
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1814

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H89").Value = 1000
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 9000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -3072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 439
$ws.Range("I97").Value = 632.3333
$ws.Range("J97").Value = 149
$ws.Range("K97").Value = 1896.9999
$ws.Range("L97").Value = 447
$ws.Range("M97").Value = -1400.9999
$ws.Range("N97").Value = -1439

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 7700
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 7700
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 23100
$ws.Range("N101").Value = -27968

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2948.1
$ws.Range("I103").Value = 1012.5
$ws.Range("J103").Value = 3432
$ws.Range("K103").Value = 3037.5
$ws.Range("L103").Value = 10296
$ws.Range("M103").Value = -2158.5
$ws.Range("N103").Value = -12054

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 4011.6
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 4011.6
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 12034.8
$ws.Range("N106").Value = -13926.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15385878
$ws.Range("I131").Value = 111111790
$ws.Range("J131").Value = 1356.4286
$ws.Range("K131").Value = 333335370
$ws.Range("L131").Value = 4069.2858
$ws.Range("M131").Value = -333330330
$ws.Range("N131").Value = -14149.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 7603.5
$ws.Range("I137").Value = 1400.8334
$ws.Range("J137").Value = 11738.611
$ws.Range("K137").Value = 4202.5002
$ws.Range("L137").Value = 35215.833
$ws.Range("M137").Value = 897.4997999999996
$ws.Range("N137").Value = -45415.833

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 198
$ws.Range("I2").Value = 243
$ws.Range("J2").Value = 159.42857
$ws.Range("K2").Value = 243
$ws.Range("L2").Value = 159.42857
$ws.Range("M2").Value = -130
$ws.Range("N2").Value = -385.42857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2920
$ws.Range("I102").Value = 3926.6667
$ws.Range("J102").Value = 2057.1428
$ws.Range("K102").Value = 3926.6667
$ws.Range("L102").Value = 2057.1428
$ws.Range("M102").Value = -2304.6667
$ws.Range("N102").Value = -5301.1428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3630.3076
$ws.Range("I122").Value = 2456.2856
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 7368.8568
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -4918.8568
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3727.75
$ws.Range("I132").Value = 3359.4443
$ws.Range("J132").Value = 4832.6665
$ws.Range("K132").Value = 10078.3329
$ws.Range("L132").Value = 14497.9995
$ws.Range("M132").Value = -7548.332900000001
$ws.Range("N132").Value = -19557.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6017.55
$ws.Range("I136").Value = 7133.8125
$ws.Range("J136").Value = 1552.5
$ws.Range("K136").Value = 21401.4375
$ws.Range("L136").Value = 4657.5
$ws.Range("M136").Value = -18851.4375
$ws.Range("N136").Value = -9757.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 44348.332
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 44348.332
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 44348.332
$ws.Range("N139").Value = -54628.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3738
$ws.Range("I81").Value = 3042.2222
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 6084.4444
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = -5023.4444
$ws.Range("N81").Value = -22122

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3738
$ws.Range("I84").Value = 3042.2222
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 30422.222
$ws.Range("L84").Value = 100000
$ws.Range("M84").Value = -25118.222
$ws.Range("N84").Value = -110608

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 52014020
$ws.Range("I122").Value = 52014020
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 156042060
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -156039610
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1798.9387
$ws.Range("I132").Value = 1412.6364
$ws.Range("J132").Value = 5198.4
$ws.Range("K132").Value = 4237.9092
$ws.Range("L132").Value = 15595.2
$ws.Range("M132").Value = -1707.9092
$ws.Range("N132").Value = -20655.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 571.7692
$ws.Range("I136").Value = 285.375
$ws.Range("J136").Value = 1030
$ws.Range("K136").Value = 856.125
$ws.Range("L136").Value = 3090
$ws.Range("M136").Value = 1693.875
